$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    $escaped = $newValue -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue "D2" '300.98'
Set-TextValue "E2" '-0.10%'
Set-TextValue "D3" '32.25'
Set-TextValue "E3" '1.83%'
Set-TextValue "D4" '4.991'
Set-TextValue "E4" '-1.97%'
Set-TextValue "D5" '0.07626'
Set-TextValue "E5" '-2.33%'
Set-TextValue "D6" '1.949'
Set-TextValue "E6" '-13.14%'
Set-TextValue "E7" '0.40%'
Set-TextValue "D8" '3.783'
Set-TextValue "E8" '-0.88%'
Set-TextValue "D9" '0.9184'
Set-TextValue "E9" '0.06%'
Set-TextValue "D10" '0.1772'
Set-TextValue "E10" '0.82%'
Set-TextValue "D11" '0.07825'
Set-TextValue "E11" '3.82%'
Set-TextValue "D12" '0.08509'
Set-TextValue "E12" '-5.39%'
Set-TextValue "D13" '0.03164'
Set-TextValue "E13" '4.42%'
Set-TextValue "D14" '0.09998'
Set-TextValue "E14" '-0.33%'
Set-TextValue "E15" '0.00%'
Set-TextValue "D16" '0.005886'
Set-TextValue "E16" '-2.69%'
Set-TextValue "E18" '-0.22%'
Set-TextValue "D19" '2.151'
Set-TextValue "E19" '-4.42%'
Set-TextValue "E21" '-2.78%'
Set-TextValue "D22" '4.281'
Set-TextValue "E22" '1.39%'
Set-TextValue "E23" '9.61%'
Set-TextValue "D24" '0.04515'
Set-TextValue "E24" '-1.50%'
Set-TextValue "D25" '0.001220'
Set-TextValue "E25" '-2.35%'
Set-TextValue "D26" '0.004387'
Set-TextValue "E26" '-1.96%'
Set-TextValue "D27" '0.0001250'
Set-TextValue "E27" '0.12%'
Set-TextValue "D39" '0.01701'
Set-TextValue "E39" '-3.88%'
Set-TextValue "D40" '0.04675'
Set-TextValue "E40" '-2.56%'
Set-TextValue "D41" '0.007504'
Set-TextValue "E41" '1.50%'
Set-TextValue "D42" '0.1348'
Set-TextValue "E42" '-0.93%'
Set-TextValue "D43" '0.002331'
Set-TextValue "E43" '6.52%'
Set-TextValue "D44" '0.01044'
Set-TextValue "E44" '2.12%'
Set-TextValue "D45" '0.00006254'
Set-TextValue "E45" '0.35%'
Set-TextValue "D46" '0.00000000750'
Set-TextValue "E46" '0.13%'
Set-TextValue "B47" 'CoinbaseStockToken'
Set-TextValue "C47" 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue "D47" '0.003105'
Set-TextValue "E47" '-61.13%'
Set-TextValue "B48" 'BOLO'
Set-TextValue "C48" 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue "D48" '0.8204'
Set-TextValue "E48" '10.84%'
Set-TextValue "D49" '0.00002101'
Set-TextValue "E49" '0.13%'
Set-TextValue "D50" '0.0002001'
Set-TextValue "E50" '0.13%'
